$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 141, shifting existing rows 141:263 down to 142:264.
# Excel's default Insert behavior copies formatting (e.g. the date number format
# on column D) from the row above, matching the original workbook's styling.
$ws.Rows(141).Insert()

# Populate the newly inserted row 141 with the new weekly price record.
$ws.Range("A141").Value = 3
$ws.Range("B141").Value = "Femacal de La Calera"
$ws.Range("C141").Value = "Coquimbo"
$ws.Range("D141").Value = 44566
$ws.Range("E141").Value = 5
$ws.Range("F141").Value = 100112009
$ws.Range("G141").Value = "Acelga"
$ws.Range("H141").Value = "Sin especificar"
$ws.Range("I141").Value = "Primera"
$ws.Range("J141").Value = 340
$ws.Range("K141").Value = 2300
$ws.Range("L141").Value = 2500
$ws.Range("M141").Value = 2406
$ws.Range("N141").Value = "`$/docena de atados (6 kilos)"
$ws.Range("O141").Value = "Provincia de Quillota"
$ws.Range("P141").Value = 401
$ws.Range("Q141").Value = 6
$ws.Range("R141").Value = "Hortaliza"
